$wb = $excel.ActiveWorkbook

# Rename "Sheet3" to "Mdm数据库脚本"
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Name = "Mdm数据库脚本"

# Make it the active sheet, scroll so row 62 is the topmost visible row,
# and select L84 (matching the saved view state of the sheet).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 62
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L84").Select()
